$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.155.93"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.786.50"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.547"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.84"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0689"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "2.044.54"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.05"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "1.792.91"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").Value = "34.116.35"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.42"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.05%  "
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.46"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.32"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0520"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "1.447.51"
$ws.Range("E35").Value = "  +5.00%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.654"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.44"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +10.49%  "
$ws.Range("E38").Value = "  +4.33%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.56"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.85%  "
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.57"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").Value = "  +4.48%  "
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "0.0₆0135"
$ws.Range("E48").Value = "  -3.42%  "
$ws.Range("D49").Value = "1.946.21"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.85"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("E51").Value = "  +0.06%  "
